{"js": "// Update the date line and the 25 division-problem answers in the table.\n// Each entry below is keyed by its OLD text so the replacement is\n// unambiguous (several new values equal other cells' old values, so a\n// blind global find/replace-all could cascade incorrectly). We resolve\n// each target by (table row, column) position instead, which exactly\n// matches the structure in the source document.\n\nconst body = context.document.body;\n\n// --- 1. Date paragraph (first paragraph in the body, outside the table) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text !== \"2025-12-17 Wednesday\") {\n  dateParagraph.insertText(\"2025-12-17 Wednesday\", \"Replace\");\n}\n\n// --- 2. Table of division answers ---\n// Raw table rows (0-based) that actually hold text; the rows in-between\n// are blank spacer rows. Each data row has 5 columns.\nconst table = body.tables.items[0];\n\nconst newValues = [\n  [\"23\u00f77=3, 2\", \"47\u00f72=23, 1\", \"56\u00f78=7, 0\", \"66\u00f79=7, 3\", \"80\u00f72=40, 0\"],\n  [\"21\u00f74=5, 1\", \"44\u00f75=8, 4\", \"28\u00f78=3, 4\", \"71\u00f79=7, 8\", \"65\u00f79=7, 2\"],\n  [\"27\u00f75=5, 2\", \"59\u00f72=29, 1\", \"71\u00f78=8, 7\", \"88\u00f76=14, 4\", \"17\u00f79=1, 8\"],\n  [\"97\u00f76=16, 1\", \"28\u00f78=3, 4\", \"71\u00f72=35, 1\", \"78\u00f77=11, 1\", \"67\u00f75=13, 2\"],\n  [\"32\u00f78=4, 0\", \"67\u00f73=22, 1\", \"50\u00f75=10, 0\", \"71\u00f79=7, 8\", \"76\u00f74=19, 0\"],\n];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const rowIdx = dataRowIndexes[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(rowIdx, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem answers in the table.\n# Targets are resolved by (row, column) position rather than a blind\n# find/replace-all, because several NEW values equal other cells' OLD\n# values (a naive global replace could cascade and corrupt later cells).\n\n$d = $word.ActiveDocument\n\n# --- 1. Date paragraph (first paragraph in the body, outside the table) ---\n$d.Paragraphs(1).Range.Text = \"2025-12-17 Wednesday\"\n\n# --- 2. Table of division answers ---\n# Raw table rows (1-based) that actually hold text; the rows in-between\n# are blank spacer rows. Each data row has 5 columns.\n$table = $d.Tables(1)\n\n$newValues = @(\n    @(\"23\u00f77=3, 2\", \"47\u00f72=23, 1\", \"56\u00f78=7, 0\", \"66\u00f79=7, 3\", \"80\u00f72=40, 0\"),\n    @(\"21\u00f74=5, 1\", \"44\u00f75=8, 4\", \"28\u00f78=3, 4\", \"71\u00f79=7, 8\", \"65\u00f79=7, 2\"),\n    @(\"27\u00f75=5, 2\", \"59\u00f72=29, 1\", \"71\u00f78=8, 7\", \"88\u00f76=14, 4\", \"17\u00f79=1, 8\"),\n    @(\"97\u00f76=16, 1\", \"28\u00f78=3, 4\", \"71\u00f72=35, 1\", \"78\u00f77=11, 1\", \"67\u00f75=13, 2\"),\n    @(\"32\u00f78=4, 0\", \"67\u00f73=22, 1\", \"50\u00f75=10, 0\", \"71\u00f79=7, 8\", \"76\u00f74=19, 0\")\n)\n$dataRowIndexes = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $dataRowIndexes.Length; $r++) {\n    $rowIdx = $dataRowIndexes[$r]\n    for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n        $table.Cell($rowIdx, $c + 1).Range.Text = $newValues[$r][$c]\n    }\n}\n"}
